# Add 2022-Q3 data:
#  1. Insert a brand-new worksheet "2022-Q3" right after "总计" and fill it
#     with the quarterly fund-holding breakdown.
#  2. Update the "总计" (summary) sheet: add a new top data row for 2022-Q3
#     and shift the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet right after "总计"
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $zj)
$q3.Name = "2022-Q3"

# Header row
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$headerRange = $q3.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows (基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名)
$codes  = @("010387","010388","005805","011157","007718","470888","011158","011453","008861","010500","010093")
$names  = @("易方达医药生物股票A","易方达医药生物股票C","华泰柏瑞医疗健康混合A","弘毅远方港股通智选领航混合A","中银创新医疗混合A","汇添富香港优势精选混合（QDII）","弘毅远方港股通智选领航混合C","华泰柏瑞医疗健康混合C","西部利得港股通新机遇灵活配置混合A","中银创新医疗混合C","西部利得港股通新机遇灵活配置混合C")
$sizes  = @("16.27","4.82","4.09","2.55","2.82","1.63","0.75","0.43","0.23","0.18","0.11")
$totpos = @("84.59","84.59","93.86","88.63","92.35","78.50","88.63","93.86","74.85","92.35","74.85")
$pospct = @("4.09","4.09","4.77","4.49","4.02","4.67","4.49","4.77","3.46","4.02","3.46")
$mktval = @("0.6654","0.1971","0.1951","0.1145","0.1134","0.0761","0.0337","0.0205","0.0080","0.0072","0.0038")
$ranks  = @(7,7,2,2,5,5,2,2,7,5,7)

for ($i = 0; $i -lt $codes.Length; $i++) {
    $r = $i + 2

    $idxCell = $q3.Cells.Item($r, 1)
    $idxCell.Value = $i
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    $codeCell = $q3.Cells.Item($r, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $codes[$i]

    $q3.Cells.Item($r, 3).Value = $names[$i]

    $sizeCell = $q3.Cells.Item($r, 4)
    $sizeCell.NumberFormat = "@"
    $sizeCell.Value = $sizes[$i]

    $totposCell = $q3.Cells.Item($r, 5)
    $totposCell.NumberFormat = "@"
    $totposCell.Value = $totpos[$i]

    $pospctCell = $q3.Cells.Item($r, 6)
    $pospctCell.NumberFormat = "@"
    $pospctCell.Value = $pospct[$i]

    $mktvalCell = $q3.Cells.Item($r, 7)
    $mktvalCell.NumberFormat = "@"
    $mktvalCell.Value = $mktval[$i]

    $q3.Cells.Item($r, 8).Value = $ranks[$i]
}

$q3.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. Update "总计" sheet: insert 2022-Q3 total as the new first data row,
#    everything else shifts down one row (re-enter explicitly to avoid
#    relying on a row-insert/shift operation).
# ---------------------------------------------------------------------
$labels = @("2022-Q3","2022-Q2","2022-Q1","2021-Q4","2021-Q3","2021-Q2","2021-Q1")
$counts = @(11,8,2,2,6,4,7)
$values = @(1.43,1.69,1.07,1.77,4.36,2.99,3.09)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = $i + 2
    $zj.Cells.Item($r, 1).Value = $i
    $zj.Cells.Item($r, 2).Value = $labels[$i]
    $zj.Cells.Item($r, 3).Value = $counts[$i]
    $zj.Cells.Item($r, 4).Value = $values[$i]
}

$zj.Range("A1").Select()
